$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.68%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.39%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.767"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.66%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06288"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.97%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.920"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.33%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.311"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "37.70%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8703"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.34%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1561"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.54%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05008"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.01%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07465"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.81%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02911"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.21%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09049"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.15%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001564"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.48%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006371"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.55%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005843"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.35%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.448"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.43%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.313"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.60%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.14%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.73%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1319"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.42%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.920"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.40%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04384"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.81%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.55%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003708"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-13.48%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001618"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-4.13%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04070"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.53%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007052"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "6.09%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1168"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.19%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002031"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.55%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01123"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-9.69%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005179"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.85%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.51%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.486"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.51%"
